$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "hskumbhar@pict.edu"
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "2025-08-29 14:30:41"
